$d = $word.ActiveDocument

# 1. Exploratory data analysis intro paragraph
$d.Content.Find.Execute("In order to build the predective model, the", $true, $false, $false, $false, $false, $true, 1, $false, "This section describes how the data has been loaded into R data.frame structures and, from the raw database, the process for extracting the features to build the predictive models.", 2)

# 2. "for the testing sets" -> "for the test sets."
$d.Content.Find.Execute("for the testing sets", $true, $false, $false, $false, $false, $true, 1, $false, "for the test sets.", 2)

# 3. Insert new runs after "data points." in the loading-the-data section:
#    " testset contains 20 unknown test cases to be predicted by the model."
$rng = $d.Content
$rng.Find.Execute("data points.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" ")
$rng.Collapse(0)
$rng.InsertAfter("testset")
$rng.Style = "Verbatim Char"
$rng.Collapse(0)
$rng.InsertAfter(" ")
$rng.Collapse(0)
$rng.InsertAfter("contains 20 unknown test cases to be predicted by the model.")

# 4. "The first step is to check ... I created a" -> "... has been to check ... A"
$d.Content.Find.Execute("The first step is to check what is the percentage of available data for each feature as the data.frame columns may contain not valid elements. I created a", $true, $false, $false, $false, $false, $true, 1, $false, "The first step has been to check what is the percentage of available data for each feature as the data.frame columns may contain not valid elements. A", 2)

# 5. "function that," -> "function has been created that,"
$d.Content.Find.Execute("function that,", $true, $false, $false, $false, $false, $true, 1, $false, "function has been created that,", 2)

# 6. "). I run the function on the" -> "). The function has been used on the"
$d.Content.Find.Execute("). I run the function on the", $true, $false, $false, $false, $false, $true, 1, $false, "). The function has been used on the", 2)

# 7. "The plot shows that data.frame variables contain or" -> add "(features)"
$d.Content.Find.Execute("The plot shows that data.frame variables contain or", $true, $false, $false, $false, $false, $true, 1, $false, "The plot shows that data.frame variables (features) contain or", 2)

# 8. Change the math run "5/" to "5%" (the discharge threshold) in the same paragraph
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $para = $paras.Item($i)
    if ($para.Range.Text -like "*features containing less the*") {
        $omaths = $para.Range.OMaths
        for ($j = 1; $j -le $omaths.Count; $j++) {
            $om = $omaths.Item($j)
            if ($om.Range.Text -eq "5/") {
                $om.Range.Text = "5%"
            }
        }
    }
}

# 9. "Moreover, we can further exclude..." paragraph rewrite
$d.Content.Find.Execute("Moreover, we can further exclude the first 7 features as they containg temporal information that has been chosen not to consider as the analysis is not considering a forcastin approach (that would be interesting to study further but it’s out of scope the present project).", $true, $false, $false, $false, $false, $true, 1, $false, "Moreover, the first 7 features contain temporal information that has been chosen not to be considered in this project (a forcasting approach would be more suitable).", 2)

# 10. "As we can see, the dataset dimension..." paragraph rewrite
$d.Content.Find.Execute("As we can see, the dataset dimension, and so its complexity, has been reduced make it also more parsimonious in its analisys.", $true, $false, $false, $false, $false, $true, 1, $false, "As we can see, the dataset dimension, and so its complexity, has been reduced making also more parsimonious its analisys.", 2)

# 11. "We split the" -> "The"
$d.Content.Find.Execute("We split the", $true, $false, $false, $false, $false, $true, 1, $false, "The", 2)

# 12. "in two," -> "has been splitted in two,"
$d.Content.Find.Execute("in two,", $true, $false, $false, $false, $false, $true, 1, $false, "has been splitted in two,", 2)

# 13. "of which will be used to train the different models and" -> "to train the different models and"
$d.Content.Find.Execute("of which will be used to train the different models and", $true, $false, $false, $false, $false, $true, 1, $false, "to train the different models and", 2)

# 14. Cross Validation paragraph rewrite
$d.Content.Find.Execute("As we will compare different algorithms, a preset Cross Validation parameter is set for all different models. A basic cross validation choise for this kind of dataset is 5-fold cross-validation to estimate accuracy. In order to seek a better estimate, each algorithm will be repeated 3 times.", $true, $false, $false, $false, $false, $true, 1, $false, "As we will compare different algorithms, a preset Cross Validation parameter is set for all different models. Since the training dataset contains a sufficient number of points, a basic cross validation choise for this kind of dataset is 5-fold cross-validation to estimate accuracy. In order to seek a better estimate, each algorithm will be repeated 3 times on each folder.", 2)

# 15. "References and Future work" -> "References and Future Work"
$d.Content.Find.Execute("References and Future work", $true, $false, $false, $false, $false, $true, 1, $false, "References and Future Work", 2)
